$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("保險")

# --- Header row (row 1): new columns reflect the dataset-standard layout ---
$ws.Range("B1").Value = "company"
$ws.Range("C1").Value = "name"
$ws.Range("D1").Value = "owner"
$ws.Range("E1").Value = "property_category"
$ws.Range("F1").Value = "category"
$ws.Range("G1").Value = "date"
$ws.Range("H1").Value = "legislator_name"
$ws.Range("I1").Value = "legislator_id"
$ws.Range("J1").Value = "source_file"
$ws.Range("K1").Value = "index"

# --- Data rows 2-17 ---
$rows = @(
    @{ row=2; A=137; B="南山人壽g"; C="南山伴我一生變額壽險"; D="洪佳君"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=137 },
    @{ row=3; A=138; B="南山人壽"; C="南山好吉利21年期還本養老保險"; D="洪佳君"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=138 },
    @{ row=4; A=139; B="南山人壽"; C="南山新新增額養老保險"; D="洪佳君"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=139 },
    @{ row=5; A=140; B="南山人壽"; C="南山好吉利21年期還本養老保險"; D="洪佳君"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=140 },
    @{ row=6; A=141; B="中泰人壽"; C="中泰人壽金富貴外幣變額年金保險"; D="洪佳君"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=141 },
    @{ row=7; A=142; B="安聯人壽"; C="世界觀外幣變額萬能壽險"; D="洪佳君"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=142 },
    @{ row=8; A=143; B="國泰人壽"; C="國泰美滿人生312終身壽險"; D="洪佳君"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=143 },
    @{ row=9; A=144; B="中華郵政"; C="六年期吉利保險"; D="洪佳君"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=144 },
    @{ row=10; A=145; B="中華郵政"; C="六年期吉利保險"; D="黃志雄"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=145 },
    @{ row=11; A=146; B="中華郵政"; C="六年期吉利保險"; D="黃志雄"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=146 },
    @{ row=12; A=147; B="新光人壽"; C="新光人壽全心終身還本保險"; D="黃志雄"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=147 },
    @{ row=13; A=148; B="新光人壽"; C="新光人壽全意終身還本保險"; D="黃志雄"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=148 },
    @{ row=14; A=149; B="南山人壽"; C="南山人壽鴻利發還本終身分紅保險"; D="黃志雄"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=149 },
    @{ row=15; A=150; B="新光人壽"; C="新光人壽全心終身還本保險"; D="黃〇誼"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=150 },
    @{ row=16; A=151; B="新光人壽"; C="新光人壽全意終身還本保險"; D="黃〇誼"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=151 },
    @{ row=17; A=152; B="富邦人壽"; C="富邦人壽安心一生终身醫療保險"; D="洪佳君"; E="insurance"; F="normal"; G="2013-12-30"; H="黃志雄"; I=1366; J="tmp22e71"; K=152 }
)

foreach ($item in $rows) {
    $r = $item.row
    $ws.Cells.Item($r, 1).Value = $item.A
    $ws.Cells.Item($r, 2).Value = $item.B
    $ws.Cells.Item($r, 3).Value = $item.C
    $ws.Cells.Item($r, 4).Value = $item.D
    $ws.Cells.Item($r, 5).Value = $item.E
    $ws.Cells.Item($r, 6).Value = $item.F
    $ws.Cells.Item($r, 7).Value = $item.G
    $ws.Cells.Item($r, 8).Value = $item.H
    $ws.Cells.Item($r, 9).Value = $item.I
    $ws.Cells.Item($r, 10).Value = $item.J
    $ws.Cells.Item($r, 11).Value = $item.K
}
